$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base value in E4; dependent formula cells (E5:E12, H4:H12) will recalc automatically.
$ws.Range("E4").Value = 4

# Update the active selection to match the authored state.
$ws.Range("E5").Select()
